# Update the "想去人数" (number of people interested) values in column F
# for the "展览" (sheet1) and "全部类型" (sheet4) worksheets.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> { row -> new F value }
$updates = @{
    "展览" = @{
        2  = 13700
        3  = 88
        5  = 541
        6  = 497
        8  = 1022
        9  = 13879
        10 = 14708
        19 = 58
        21 = 1141
        24 = 5683
        25 = 942
        27 = 5393
        29 = 44
        30 = 237
    }
    "全部类型" = @{
        2  = 13700
        3  = 88
        6  = 541
        7  = 497
        9  = 1022
        10 = 13879
        11 = 14708
        20 = 58
        22 = 1141
        25 = 5683
        26 = 942
        28 = 5393
        30 = 44
        31 = 237
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowsForSheet = $updates[$sheetName]
    foreach ($row in $rowsForSheet.Keys) {
        $newValue = $rowsForSheet[$row]
        $ws.Range("F$row").Value = $newValue
    }
}
